# Apply the edit described in the diff.
#
# Both worksheets ("backtracking" and "sat") have an "ADDITIONAL CASES" table
# (rows 32-40) with one row per puzzle size: 6x6, 7x7, 8x8, 9x9, 10x10, 11x11,
# 12x12. The edit removes the "6x6" row and the "12x12" row, so the remaining
# rows (7x7 .. 11x11) shift up and the table shrinks from rows 34-40 down to
# rows 34-38.
#
# It also switches which sheet/window state is active: the "sat" sheet becomes
# the selected/active tab (instead of "backtracking"), and each sheet's
# remembered selection changes.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("backtracking")
$ws2 = $wb.Worksheets.Item("sat")

foreach ($ws in @($ws1, $ws2)) {
    # Delete the "6x6" row (row 34). Everything below (7x7 .. 12x12) shifts up
    # one row, so the table now occupies rows 34-39 with "12x12" now on row 39.
    $ws.Rows.Item(34).Delete()

    # Delete the (now shifted) "12x12" row, leaving 7x7 .. 11x11 on rows 34-38.
    $ws.Rows.Item(39).Delete()

    # Re-fill the average column as a single range formula so the engine
    # re-groups E34:E38 into one shared formula (si/ref), matching the way
    # Excel keeps a shared AVERAGE formula intact across the shifted rows.
    $ws.Range("E34:E38").Formula = "=AVERAGE(B34:D34)"
}

# "sat" is now the active/selected sheet (was "backtracking" before).
$ws2.Select()
$ws2.Activate()

# Each sheet keeps its own remembered selection. "sat" is selected last so
# the workbook's active-cell state lines up with the now-active "sat" sheet.
$ws1.Range("G35").Select()
$ws2.Range("L5").Select()
